{"js": "// Add a \"Monsters\" paragraph plus a linked Roll20 monster-compendium\n// paragraph right after the existing StarterSet_Characters.pdf hyperlink\n// paragraph (matches the commit \"added npcs and some scene descriptions\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the paragraph that hosts the StarterSet_Characters.pdf hyperlink \u2014\n// the new content is inserted directly after it.\nlet anchorParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.indexOf(\"https://media.wizards.com/downloads/dnd/StarterSet_Characters.pdf\") !== -1) {\n    anchorParagraph = paragraph;\n    break;\n  }\n}\nif (!anchorParagraph) {\n  // Fall back to the last paragraph if the expected link text isn't found.\n  anchorParagraph = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// New plain-text paragraph: \"Monsters\"\nconst monstersParagraph = anchorParagraph.insertParagraph(\"Monsters\", \"After\");\n\n// New paragraph containing the Roll20 monsters hyperlink.\nconst monsterLinkUrl = \"https://roll20.net/compendium/dnd5e/Monsters#content\";\nconst linkParagraph = monstersParagraph.insertParagraph(monsterLinkUrl, \"After\");\nlinkParagraph.getRange().hyperlink = monsterLinkUrl;\n\nawait context.sync();\n", "ps1": "# Add a \"Monsters\" paragraph plus a linked Roll20 monster-compendium\n# paragraph right after the existing StarterSet_Characters.pdf hyperlink\n# paragraph (matches the commit \"added npcs and some scene descriptions\").\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that hosts the StarterSet_Characters.pdf hyperlink \u2014\n# the new content is inserted directly after it.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*https://media.wizards.com/downloads/dnd/StarterSet_Characters.pdf*\") {\n        $target = $p\n    }\n}\nif ($target -eq $null) {\n    # Fall back to the document's last paragraph if the expected link text\n    # isn't found.\n    $target = $d.Paragraphs.Last\n}\n\n# New plain-text paragraph: \"Monsters\"\n$targetEnd = $target.Range\n$targetEnd.Collapse(0)   # wdCollapseEnd\n$targetEnd.InsertParagraphAfter()\n$monstersPara = $target.Next()\n$monstersPara.Range.Text = \"Monsters\"\n\n# New paragraph containing the Roll20 monsters hyperlink.\n$monsterLinkUrl = \"https://roll20.net/compendium/dnd5e/Monsters#content\"\n\n$monstersEnd = $monstersPara.Range\n$monstersEnd.Collapse(0)  # wdCollapseEnd\n$monstersEnd.InsertParagraphAfter()\n$linkPara = $monstersPara.Next()\n\n$linkPara.Range.InsertAfter($monsterLinkUrl)\n$linkRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start + $monsterLinkUrl.Length)\n$d.Hyperlinks.Add($linkRange, $monsterLinkUrl, $null, $null, $null, $null) | Out-Null\n"}
